# Applies the "commits as of 6thmay2020" change:
#  - Adds two new worksheets: "RecipientQuery" and "AddressBookQuery",
#    each holding a "Query" header cell and a wrapped multi-line SQL
#    string underneath, with a widened/tall cell to show the text.
#  - As a natural side effect of adding sheets (and leaving the last
#    one active/selected), the workbook's active-tab / tabSelected
#    bookkeeping moves from the first sheet to the newly added last
#    sheet, matching the diff's sheet1 tabSelected removal + new-sheet
#    tabSelected addition and bookViews activeTab update.

$wb = $excel.ActiveWorkbook

$nl = [char]10

# ---- New sheet 1: RecipientQuery ----------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$recipientQuery = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$recipientQuery.Name = "RecipientQuery"

$recipientQuery.Range("A1").Value = "Query"

$recipientSql = "SELECT " + $nl + "    [FirstName] as 'Name1'" + $nl + "    ,[LastName] as 'Name2'" + $nl + "    ,[FaxNumber] as 'Fax Number'" + $nl + "    ,[LastChangedBy] as 'Last Changed By'" + $nl + "    ,[LastChangedOn] as 'Last Changed On'" + $nl + "    FROM [Fax_Recipient]"
$recipientQuery.Range("A2").Value = $recipientSql
$recipientQuery.Range("A2").WrapText = $true

$recipientQuery.Columns.Item(1).ColumnWidth = 55.17
$recipientQuery.Rows.Item(2).RowHeight = 105

$recipientQuery.Range("A18").Select() | Out-Null

# ---- New sheet 2: AddressBookQuery --------------------------------
$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$addressBookQuery = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet2)
$addressBookQuery.Name = "AddressBookQuery"

$addressBookQuery.Range("A1").Value = "Query"

$addressSql = "SELECT " + $nl + "  [Name] as Name" + $nl + "  ,[FaxLine] as 'Fax Line'" + $nl + "  ,[LastChangedBy] as 'Last Changed By'" + $nl + "  ,[LastChangedOn] as 'Last Changed On'" + $nl + "  FROM [Fax_AddressBook]"
$addressBookQuery.Range("A2").Value = $addressSql
$addressBookQuery.Range("A2").WrapText = $true

$addressBookQuery.Columns.Item(1).ColumnWidth = 37.8
$addressBookQuery.Rows.Item(2).RowHeight = 90

$addressBookQuery.Range("A10").Select() | Out-Null
